$wb = $excel.ActiveWorkbook
$main = $wb.Worksheets.Item("Main")
$model = $wb.Worksheets.Item("Model")

# --- Main sheet: new commentary / notes in column H, competitor & M&A lists ---
# (cells are populated in the same order the shared-string table records them)

# Company description blurb (row 4)
$main.Range("H4").Value = "NVIDIA is a semiconductor company based in Santa Clara, California. NVIDIA has made significant advancements in AI and DL through it's CUDA platform."

# Jensen Huang notes block (rows 11-16)
$main.Range("H11").Value = "Founder and CEO of NVIDIA since inception"
$main.Range("H12").Value = "Bachelors Electrical Eng @ Oregon State"
$main.Range("H13").Value = "Masters EE Stanford"
$main.Range("H14").Value = "60 yrs old"
$main.Range("H15").Value = "owns 3.6% of company stock"
$main.Range("H16").Value = "21.36M yearly comp"

# Competitors list (H19:H25)
$main.Range("H19").Value = "Reneas"
$main.Range("H20").Value = "Cisco"
$main.Range("H21").Value = "AMD"
$main.Range("H22").Value = "HP"
$main.Range("H23").Value = "Broadcom"
$main.Range("H24").Value = "Intel"
$main.Range("H25").Value = "Qualcomm"

# "Competitors:" header moves from G16 down to G18 (keep bold style on both cells)
$main.Range("G16").ClearContents()
$main.Range("G18").Value = "Competitors:"
$main.Range("G18").Font.Bold = $true

# New "M&A:" header at J18, plus a blank styled cell at M18 matching
$main.Range("J18").Value = "M&A:"
$main.Range("J18").Font.Bold = $true
$main.Range("M18").Font.Bold = $true

# M&A list (K19:K24)
$main.Range("K19").Value = "Mellanox "
$main.Range("K20").Value = "Arm"
$main.Range("K21").Value = "SwiftStack"
$main.Range("K22").Value = "Ageia"
$main.Range("K23").Value = "OmniML"
$main.Range("K24").Value = "CoreWeave"

# TSMC note (row 5) added last
$main.Range("H5").Value = "TSMC makes NVIDIA chips"

# --- Model sheet: discount / growth assumption update ---
$model.Range("AI6").Value = 0.08

# --- View state: Model sheet selection/pane moves to AI13, Main becomes the active tab ---
$model.Activate()
$model.Range("AI13").Select()

$main.Activate()
$main.Range("E22").Select()
